# Traffic-light phasing sheet: replace the second data row's key/signal
# values with the new real-time phasing data, drop the now-obsolete rows
# (the old rows 3-5), and re-split column B's width from the combined
# B:C formatting so it can hold the longer values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (the only row that survives) with the new phase strings.
$ws.Range("B2").Value = "srgsrrrsrgsrrr"
$ws.Range("C2").Value = "srrsrrgsrrsrrg"
$ws.Range("D2").Value = "6,6,1,8,8,8,3,2,2,5,4,4,4,7"

# Remove the old rows 3-5 entirely (their data has been superseded).
$ws.Range("A3:A5").EntireRow.Delete()

# Column B now needs its own (wider) width, split out from the old
# shared B:C formatting.
$ws.Columns.Item(2).ColumnWidth = 14

# Match the saved selection/active cell.
$ws.Range("A2").Select()
